$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.138.17'
$ws.Range("E2").Value = '  -1.88%  '
$ws.Range("D3").Value = '1.656.41'
$ws.Range("E3").Value = '  -1.89%  '
$ws.Range("E4").Value = '  +0.40%  '
$ws.Range("D5").Value = '''217.54'
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("D6").Value = '''0.5213'
$ws.Range("E6").Value = '  -2.38%  '
$ws.Range("E7").Value = '  +0.49%  '
$ws.Range("D8").Value = '''0.2659'
$ws.Range("E8").Value = '  -0.88%  '
$ws.Range("D9").Value = '''0.06296'
$ws.Range("E9").Value = '  -2.14%  '
$ws.Range("D10").Value = '''20.95'
$ws.Range("E10").Value = '  -3.40%  '
$ws.Range("D11").Value = '''0.07711'
$ws.Range("E11").Value = '  -1.05%  '
$ws.Range("D12").Value = '1.669.68'
$ws.Range("E12").Value = '  -1.14%  '
$ws.Range("D13").Value = '''4.413'
$ws.Range("E13").Value = '  -2.03%  '
$ws.Range("D14").Value = '1.886.86'
$ws.Range("D15").Value = '''0.5436'
$ws.Range("E15").Value = '  -3.34%  '
$ws.Range("D16").Value = '0.0₅8200'
$ws.Range("E16").Value = '  -2.94%  '
$ws.Range("D17").Value = '''64.68'
$ws.Range("E17").Value = '  -2.42%  '
$ws.Range("D18").Value = '26.205.46'
$ws.Range("E18").Value = '  -1.79%  '
$ws.Range("D19").Value = '''1.005'
$ws.Range("E19").Value = '  +0.35%  '
$ws.Range("D20").Value = '''4.648'
$ws.Range("E20").Value = '  -3.40%  '
$ws.Range("D21").Value = '''192.22'
$ws.Range("E21").Value = '  -1.74%  '
$ws.Range("D22").Value = '''10.11'
$ws.Range("E22").Value = '  -2.98%  '
$ws.Range("D23").Value = '''6.046'
$ws.Range("E23").Value = '  -5.34%  '
$ws.Range("D24").Value = '''1.009'
$ws.Range("E24").Value = '  +0.61%  '
$ws.Range("D25").Value = '''138.65'
$ws.Range("E25").Value = '  -3.79%  '
$ws.Range("D26").Value = '''0.1233'
$ws.Range("E26").Value = '  -4.56%  '
$ws.Range("D27").Value = '''7.180'
$ws.Range("E27").Value = '  -4.14%  '
$ws.Range("D28").Value = '''16.10'
$ws.Range("E28").Value = '  -1.17%  '
$ws.Range("E29").Value = '  -0.77%  '
$ws.Range("D30").Value = '''0.05980'
$ws.Range("E30").Value = '  -3.13%  '
$ws.Range("D31").Value = '''1.280'
$ws.Range("E31").Value = '  -0.17%  '
$ws.Range("D32").Value = '''3.584'
$ws.Range("E32").Value = '  -0.61%  '
$ws.Range("D33").Value = '''3.311'
$ws.Range("E33").Value = '  -4.71%  '
$ws.Range("D34").Value = '''1.641'
$ws.Range("E34").Value = '  -3.76%  '
$ws.Range("D35").Value = '''0.9762'
$ws.Range("E35").Value = '  -3.80%  '
$ws.Range("D36").Value = '''2.785'
$ws.Range("E36").Value = '  -0.50%  '
$ws.Range("D37").Value = '''2.414'
$ws.Range("E37").Value = '  -0.37%  '
$ws.Range("D38").Value = '''0.5898'
$ws.Range("E38").Value = '  +2.71%  '
$ws.Range("D39").Value = '''0.01583'
$ws.Range("E39").Value = '  -4.17%  '
$ws.Range("D40").Value = '''5.932'
$ws.Range("E40").Value = '  -1.54%  '
$ws.Range("D41").Value = '''0.8645'
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("E42").Value = '  +0.33%  '
$ws.Range("D43").Value = '1.034.57'
$ws.Range("E43").Value = '  -4.36%  '
$ws.Range("D44").Value = '''99.49'
$ws.Range("E44").Value = '  -1.02%  '
$ws.Range("D45").Value = '1.801.99'
$ws.Range("D46").Value = '''56.94'
$ws.Range("E46").Value = '  -1.01%  '
$ws.Range("D47").Value = '0.0₈106'
$ws.Range("E47").Value = '  -2.16%  '
$ws.Range("D48").Value = '''1.001'
$ws.Range("E48").Value = '  -0.09%  '
$ws.Range("D49").Value = '''8.049'
$ws.Range("E49").Value = '  -2.06%  '
$ws.Range("D50").Value = '''0.05181'
$ws.Range("E50").Value = '  -0.84%  '
$ws.Range("D51").Value = '''0.4230'
$ws.Range("E51").Value = '  -0.31%  '
